$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 7058.19682190226
$ws.Range("C2").Value = 6551.75777183103
$ws.Range("E2").Value = 2974.08885967773
$ws.Range("F2").Value = -3.9230570204686

$ws.Range("B3").Value = 7395.25178308535
$ws.Range("C3").Value = 6322.71174715432
$ws.Range("E3").Value = 3263.62188669975
$ws.Range("F3").Value = 248.59723474392

$ws.Range("B4").Value = 3130.84389747465
$ws.Range("C4").Value = 4723.65121545155
$ws.Range("E4").Value = 3660.82331626128
$ws.Range("F4").Value = 198.519772154701

$ws.Range("B5").Value = 3059.00476017793
$ws.Range("C5").Value = 4593.28822606136
$ws.Range("E5").Value = 3677.78226218023
$ws.Range("F5").Value = 193.794603676733

$ws.Range("B6").Value = 8811.83146577181
$ws.Range("C6").Value = 8115.63368752306
$ws.Range("E6").Value = 4537.4083276642
$ws.Range("F6").Value = 376.376750632802

$ws.Range("B7").Value = 8620.20120214289
$ws.Range("C7").Value = 7978.9153967638
$ws.Range("E7").Value = 4382.71816243105
$ws.Range("F7").Value = 364.234731633119

$ws.Range("E11").Value = 3796.87610080022
$ws.Range("F11").Value = -46.8270022449817

$ws.Range("E12").Value = 3796.87610080022
$ws.Range("F12").Value = -42.0614738050249

$ws.Range("E13").Value = 4685.46349322388
$ws.Range("F13").Value = 121.648267675137

$ws.Range("E14").Value = 4685.46349322388
$ws.Range("F14").Value = 116.53173213716

$ws.Range("E15").Value = 4685.46349322388
$ws.Range("F15").Value = 123.247531323142
